{"js": "// Load all paragraphs in the document body with their text so we can\n// locate the exact anchor paragraphs described by the diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Helper: find the index of the first paragraph (optionally starting the\n// search at `fromIndex`) whose text equals `text` exactly.\nfunction findIndex(text, fromIndex) {\n  const start = fromIndex || 0;\n  for (let i = start; i < items.length; i++) {\n    if (items[i].text === text) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + text);\n}\n\n// ---------------------------------------------------------------------\n// 1) Insert a new bullet before \"Statements:\" (the \"Kinds (Wrapped ...\" /\n//    \"Statements:\" boundary near the top of the document).\n// ---------------------------------------------------------------------\nconst idxStatements = findIndex(\"Statements:\");\nitems[idxStatements].insertParagraph(\n  \"Augment Resources with Kinds in Context. Core Model Transforms Mappings Instances / Roles: Kinds and Singleton (Resource) Class.\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Delete the \"(Class, Instance, Atribute, Value)\" bullet that used to\n//    precede \"Augmentations: Activation (Schema), ...\".\n// ---------------------------------------------------------------------\nconst idxClassInstance = findIndex(\"(Class, Instance, Atribute, Value)\");\nitems[idxClassInstance].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Insert 4 new bullets before the second \"Activation:\" bullet (the one\n//    that follows \"Inputs Occurrences / Occurring Augmentations: ...\" and\n//    the empty bullet after it).\n// ---------------------------------------------------------------------\nconst idxAugment = findIndex(\n  \"Inputs Occurrences / Occurring Augmentations: Alignment, Aggregation, Activation of raw Template Inputs. Matchings / Merge (data, schema, behaviors)\"\n);\n// \"Activation:\" occurs twice in the document; the one we need is the\n// first occurrence *after* the \"Inputs Occurrences / Occurring ...\" bullet.\nconst idxActivation2 = findIndex(\"Activation:\", idxAugment);\n\nlet anchor = items[idxActivation2];\nanchor = anchor.insertParagraph(\"Inputs:\", Word.InsertLocation.before);\nawait context.sync();\nanchor = anchor.insertParagraph(\n  \"Augment Resource with Kind in Context. Core Model Transforms Mappings Instances / Roles: Kinds and Singleton (Resource) Class.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nanchor = anchor.insertParagraph(\n  \"(Class : Transform, Instance : Kind T, Atribute : Mapping, Value : Kind U);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nanchor = anchor.insertParagraph(\n  \"Inputs Normal Forms: Dimensional, Discrete, etc. Parse Aggregations into Core Model.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nanchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4) Rewrite the \"Infer T / U Class / Metaclass ...\" bullet and insert two\n//    new bullets right after it.\n// ---------------------------------------------------------------------\nconst idxInfer = findIndex(\n  \"Infer T / U Class / Metaclass (Mappings / Transforms) Wrapped Types / Instances (Dimension / Time). Parse Instances (Subject / Object Resources) Wrapper / Wrapped Types.\"\n);\nlet inferPara = items[idxInfer];\ninferPara.insertText(\n  \"Infer T / U Resources Kinds Wrapped Types / Instances (Place / NY), Infer Transform Class / Metaclass by SPO Kinds. Parse Instances (Subject / Object Resources) Wrapper / Wrapped Types (Kinds Matching). Transform / Kinds Resolve Mapping Statement (noop, merge, add);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nlet afterInfer = inferPara.insertParagraph(\n  \"Matching in Occurrence / Occurring Direction.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nafterInfer = afterInfer.insertParagraph(\n  \"Resource::Mapping::Kind\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 5) Rewrite \"Assert: (Class / Transform, Resource, Attribute / Mapping,\n//    Resource);\" -- only the FIRST occurrence (inside \"Core Model:\"\n//    section) is changed; the later occurrence (inside \"Core Model\n//    Templates:\") stays untouched.\n// ---------------------------------------------------------------------\nconst idxAssert1 = findIndex(\n  \"Assert: (Class / Transform, Resource, Attribute / Mapping, Resource);\"\n);\nitems[idxAssert1].insertText(\n  \"Assert: (Class / Transform, Resource : Kind T, Attribute / Mapping, Resource : Kind U);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 6) Rewrite \"Occurrence / Mapping Declaration: ...\" bullet.\n// ---------------------------------------------------------------------\nconst idxOccurrenceDecl = findIndex(\n  \"Occurrence / Mapping Declaration: (Mapping / Class / Metaclass : T,\\u00A0 Resource / Instance, Transform / Occurrence / Context / Statement / Class / Metaclass : U, Resource / Instance / Role);\"\n);\nitems[idxOccurrenceDecl].insertText(\n  \"Occurrence / Mapping Declaration: (Mapping / Class / Metaclass,\\u00A0 Resource / Instance : T, Transform / Occurrence / Context / Statement / Class / Metaclass, Resource / Instance / Role : U);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 7) Rewrite \"Occurring / Transform Application: ...\" bullet.\n// ---------------------------------------------------------------------\nconst idxOccurringApp = findIndex(\n  \"Occurring / Transform Application: (Transform / Class / Metaclass : T, Resource / Instance, Mapping / Occurring / Context / Statement / Class / Metaclass : U, Resource / Instance / Role);\"\n);\nitems[idxOccurringApp].insertText(\n  \"Occurring / Transform Application: (Transform / Class / Metaclass, Resource / Instance : Kind T, Mapping / Occurring / Context / Statement / Class / Metaclass, Resource / Instance / Role : Kind U);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 8) Insert two new bullets before \"(Time, 1h, mins, 60m);\".\n// ---------------------------------------------------------------------\nconst idxTime = findIndex(\"(Time, 1h, mins, 60m);\");\nlet beforeTime = items[idxTime].insertParagraph(\n  \"Inputs Normal Forms: Dimensional, Discrete, etc. Parse Aggregations into Core Model.\",\n  Word.InsertLocation.before\n);\nawait context.sync();\nbeforeTime.insertParagraph(\n  \"(Class : Transform, Instance : Kind T, Atribute : Mapping, Value : Kind U);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 9) Insert an empty bullet before \"(Employment, anEmployment, employee,\n//    John);\", then insert \"(Employment, employee, Employee);\" right\n//    after that same bullet.\n// ---------------------------------------------------------------------\nconst idxEmploymentAnEmployee = findIndex(\n  \"(Employment, anEmployment, employee, John);\"\n);\nconst employmentAnEmployeePara = items[idxEmploymentAnEmployee];\nemploymentAnEmployeePara.insertParagraph(\"\", Word.InsertLocation.before);\nawait context.sync();\nemploymentAnEmployeePara.insertParagraph(\n  \"(Employment, employee, Employee);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 10) Replace \"(John, employment, anEmployment);\" bullet's text, delete\n//     \"(ABC Inc, employment, anEmployment);\" bullet, and insert 4 new\n//     bullets in its place.\n// ---------------------------------------------------------------------\nconst idxJohnEmployment = findIndex(\"(John, employment, anEmployment);\");\nconst idxAbcEmployment = findIndex(\"(ABC Inc, employment, anEmployment);\");\n\nlet johnEmploymentPara = items[idxJohnEmployment];\njohnEmploymentPara.insertText(\n  \"(Employment, employer, Employer);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nitems[idxAbcEmployment].delete();\nawait context.sync();\n\nlet afterEmployer = johnEmploymentPara.insertParagraph(\n  \"(John: Transform / Singleton, John, employment, anEmployment);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nafterEmployer = afterEmployer.insertParagraph(\n  \"(Employee, employment, Employment);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nafterEmployer = afterEmployer.insertParagraph(\n  \"(ABC Inc: Transform / Singleton, ABC Inc, employment, anEmployment);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nafterEmployer = afterEmployer.insertParagraph(\n  \"(Employer, employment, Employee);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 11) Replace \"(ABC Inc, employs, John);\" bullet's text and insert 2 new\n//     bullets right after it.\n// ---------------------------------------------------------------------\nconst idxAbcEmploys = findIndex(\"(ABC Inc, employs, John);\");\nlet abcEmploysPara = items[idxAbcEmploys];\nabcEmploysPara.insertText(\n  \"(Employee, employmentAt, Employer);\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nlet afterEmploysFor = abcEmploysPara.insertParagraph(\n  \"(ABC Inc, employsFor, John);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\nafterEmploysFor.insertParagraph(\n  \"(Employer, employsFor, Employee);\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Word COM interop script implementing the same edits as edit.js.\n# $d / $word / $app resolve via the harness; $d is the active document.\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# Helper: return the 1-based Paragraphs index of the first paragraph\n# (starting the search at $startIndex, default 1) whose text equals\n# $text exactly (ignoring the trailing paragraph-mark character).\n# ------------------------------------------------------------------\nfunction Find-ParaIndex($doc, $text, $startIndex) {\n    if (-not $startIndex) { $startIndex = 1 }\n    $count = $doc.Paragraphs.Count\n    for ($i = $startIndex; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        if ($t.Length -gt 0) {\n            $t = $t.Substring(0, $t.Length - 1)\n        }\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\n# Helper: replace the text of an existing paragraph (identified by its\n# 1-based index) while leaving its own paragraph mark untouched.\nfunction Set-ParaText($doc, $index, $newText) {\n    $p = $doc.Paragraphs.Item($index)\n    $r = $p.Range.Duplicate()\n    $r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude paragraph mark\n    $r.Text = $newText\n}\n\n# Helper: insert a brand-new paragraph with $newText right after the\n# paragraph at 1-based $index; returns the index of the newly created\n# paragraph.\nfunction Insert-ParaAfter($doc, $index, $newText) {\n    $p = $doc.Paragraphs.Item($index)\n    $p.Range.InsertParagraphAfter()\n    $newIndex = $index + 1\n    if ($newText -ne \"\") {\n        $newPara = $doc.Paragraphs.Item($newIndex)\n        $newPara.Range.InsertBefore($newText)\n    }\n    return $newIndex\n}\n\n# Helper: insert a brand-new paragraph with $newText right before the\n# paragraph at 1-based $index; returns the index of the newly created\n# paragraph (the original paragraph shifts to $index + 1).\nfunction Insert-ParaBefore($doc, $index, $newText) {\n    $p = $doc.Paragraphs.Item($index)\n    $r = $p.Range.Duplicate()\n    $r.Collapse(1) | Out-Null  # wdCollapseStart = 1\n    $r.InsertParagraphBefore()\n    if ($newText -ne \"\") {\n        $newPara = $doc.Paragraphs.Item($index)\n        $newPara.Range.InsertBefore($newText)\n    }\n    return $index\n}\n\n# ------------------------------------------------------------------\n# 1) Insert a new bullet before \"Statements:\".\n# ------------------------------------------------------------------\n$idxStatements = Find-ParaIndex $d \"Statements:\" 1\nInsert-ParaBefore $d $idxStatements \"Augment Resources with Kinds in Context. Core Model Transforms Mappings Instances / Roles: Kinds and Singleton (Resource) Class.\" | Out-Null\n\n# ------------------------------------------------------------------\n# 2) Delete the \"(Class, Instance, Atribute, Value)\" bullet.\n# ------------------------------------------------------------------\n$idxClassInstance = Find-ParaIndex $d \"(Class, Instance, Atribute, Value)\" 1\n$d.Paragraphs.Item($idxClassInstance).Range.Delete() | Out-Null\n\n# ------------------------------------------------------------------\n# 3) Insert 4 new bullets before the second \"Activation:\" bullet (the\n#    one following \"Inputs Occurrences / Occurring Augmentations: ...\").\n# ------------------------------------------------------------------\n$idxAugment = Find-ParaIndex $d \"Inputs Occurrences / Occurring Augmentations: Alignment, Aggregation, Activation of raw Template Inputs. Matchings / Merge (data, schema, behaviors)\" 1\n$idxActivation2 = Find-ParaIndex $d \"Activation:\" $idxAugment\n\n$cursor = Insert-ParaBefore $d $idxActivation2 \"Inputs:\"\n$cursor = Insert-ParaAfter $d $cursor \"Augment Resource with Kind in Context. Core Model Transforms Mappings Instances / Roles: Kinds and Singleton (Resource) Class.\"\n$cursor = Insert-ParaAfter $d $cursor \"(Class : Transform, Instance : Kind T, Atribute : Mapping, Value : Kind U);\"\n$cursor = Insert-ParaAfter $d $cursor \"Inputs Normal Forms: Dimensional, Discrete, etc. Parse Aggregations into Core Model.\"\n$cursor = Insert-ParaAfter $d $cursor \"\"\n\n# ------------------------------------------------------------------\n# 4) Rewrite the \"Infer T / U Class / Metaclass ...\" bullet and insert\n#    two new bullets right after it.\n# ------------------------------------------------------------------\n$idxInfer = Find-ParaIndex $d \"Infer T / U Class / Metaclass (Mappings / Transforms) Wrapped Types / Instances (Dimension / Time). Parse Instances (Subject / Object Resources) Wrapper / Wrapped Types.\" 1\nSet-ParaText $d $idxInfer \"Infer T / U Resources Kinds Wrapped Types / Instances (Place / NY), Infer Transform Class / Metaclass by SPO Kinds. Parse Instances (Subject / Object Resources) Wrapper / Wrapped Types (Kinds Matching). Transform / Kinds Resolve Mapping Statement (noop, merge, add);\"\n\n$cursor = Insert-ParaAfter $d $idxInfer \"Matching in Occurrence / Occurring Direction.\"\n$cursor = Insert-ParaAfter $d $cursor \"Resource::Mapping::Kind\"\n\n# ------------------------------------------------------------------\n# 5) Rewrite \"Assert: (Class / Transform, Resource, Attribute / Mapping,\n#    Resource);\" -- only the FIRST occurrence is changed.\n# ------------------------------------------------------------------\n$idxAssert1 = Find-ParaIndex $d \"Assert: (Class / Transform, Resource, Attribute / Mapping, Resource);\" 1\nSet-ParaText $d $idxAssert1 \"Assert: (Class / Transform, Resource : Kind T, Attribute / Mapping, Resource : Kind U);\"\n\n# ------------------------------------------------------------------\n# 6) Rewrite \"Occurrence / Mapping Declaration: ...\" bullet. Note: the\n#    text contains a NBSP (U+00A0) immediately after the comma that\n#    follows \"Metaclass\" -- preserved in both the old and new text.\n# ------------------------------------------------------------------\n$nbsp = [char]0x00A0\n$declOld = \"Occurrence / Mapping Declaration: (Mapping / Class / Metaclass : T,$nbsp Resource / Instance, Transform / Occurrence / Context / Statement / Class / Metaclass : U, Resource / Instance / Role);\"\n$idxOccurrenceDecl = Find-ParaIndex $d $declOld 1\n$declNew = \"Occurrence / Mapping Declaration: (Mapping / Class / Metaclass,$nbsp Resource / Instance : T, Transform / Occurrence / Context / Statement / Class / Metaclass, Resource / Instance / Role : U);\"\nSet-ParaText $d $idxOccurrenceDecl $declNew\n\n# ------------------------------------------------------------------\n# 7) Rewrite \"Occurring / Transform Application: ...\" bullet.\n# ------------------------------------------------------------------\n$idxOccurringApp = Find-ParaIndex $d \"Occurring / Transform Application: (Transform / Class / Metaclass : T, Resource / Instance, Mapping / Occurring / Context / Statement / Class / Metaclass : U, Resource / Instance / Role);\" 1\nSet-ParaText $d $idxOccurringApp \"Occurring / Transform Application: (Transform / Class / Metaclass, Resource / Instance : Kind T, Mapping / Occurring / Context / Statement / Class / Metaclass, Resource / Instance / Role : Kind U);\"\n\n# ------------------------------------------------------------------\n# 8) Insert two new bullets before \"(Time, 1h, mins, 60m);\".\n# ------------------------------------------------------------------\n$idxTime = Find-ParaIndex $d \"(Time, 1h, mins, 60m);\" 1\n$cursor = Insert-ParaBefore $d $idxTime \"Inputs Normal Forms: Dimensional, Discrete, etc. Parse Aggregations into Core Model.\"\n$cursor = Insert-ParaAfter $d $cursor \"(Class : Transform, Instance : Kind T, Atribute : Mapping, Value : Kind U);\"\n\n# ------------------------------------------------------------------\n# 9) Insert an empty bullet before \"(Employment, anEmployment, employee,\n#    John);\", then insert \"(Employment, employee, Employee);\" right\n#    after that same bullet.\n# ------------------------------------------------------------------\n$idxEmploymentAnEmployee = Find-ParaIndex $d \"(Employment, anEmployment, employee, John);\" 1\n$idxEmploymentAnEmployee = Insert-ParaBefore $d $idxEmploymentAnEmployee \"\"\n$idxEmploymentAnEmployee = $idxEmploymentAnEmployee + 1   # original paragraph shifted down by one\nInsert-ParaAfter $d $idxEmploymentAnEmployee \"(Employment, employee, Employee);\" | Out-Null\n\n# ------------------------------------------------------------------\n# 10) Replace \"(John, employment, anEmployment);\" bullet's text, delete\n#     \"(ABC Inc, employment, anEmployment);\" bullet, and insert 4 new\n#     bullets in its place.\n# ------------------------------------------------------------------\n$idxJohnEmployment = Find-ParaIndex $d \"(John, employment, anEmployment);\" 1\n$idxAbcEmployment = Find-ParaIndex $d \"(ABC Inc, employment, anEmployment);\" 1\n\nSet-ParaText $d $idxJohnEmployment \"(Employment, employer, Employer);\"\n$d.Paragraphs.Item($idxAbcEmployment).Range.Delete() | Out-Null\n\n$cursor = Insert-ParaAfter $d $idxJohnEmployment \"(John: Transform / Singleton, John, employment, anEmployment);\"\n$cursor = Insert-ParaAfter $d $cursor \"(Employee, employment, Employment);\"\n$cursor = Insert-ParaAfter $d $cursor \"(ABC Inc: Transform / Singleton, ABC Inc, employment, anEmployment);\"\n$cursor = Insert-ParaAfter $d $cursor \"(Employer, employment, Employee);\"\n\n# ------------------------------------------------------------------\n# 11) Replace \"(ABC Inc, employs, John);\" bullet's text and insert 2 new\n#     bullets right after it.\n# ------------------------------------------------------------------\n$idxAbcEmploys = Find-ParaIndex $d \"(ABC Inc, employs, John);\" 1\nSet-ParaText $d $idxAbcEmploys \"(Employee, employmentAt, Employer);\"\n\n$cursor = Insert-ParaAfter $d $idxAbcEmploys \"(ABC Inc, employsFor, John);\"\n$cursor = Insert-ParaAfter $d $cursor \"(Employer, employsFor, Employee);\"\n"}
